$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Marking" row correct-answer marks value (B11: 3 -> 5)
$ws.Range("B11").Value = 5

# Update the "Total" row marks obtained (B12: 66 -> 110)
$ws.Range("B12").Value = 110

# Update the "corr/total" summary text (E12: "64/84" -> "110/140")
$ws.Range("E12").Value = "110/140"
